$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Menn - statistikk": add two new top-result rows for Marcel Grohmann
# (Vestby IL) - a new "Kule" (shot put) result and a new "Diskos" (discus)
# result, which pushes existing rows further down the ranked list.
# ---------------------------------------------------------------------------
$wsMenn = $wb.Worksheets.Item("Menn - statistikk")

# Insert a blank row before the (old) row 19 - this is where the new "Diskos"
# entry belongs (between the "Høyde u.t" row and the "Samhald IL" row).
$wsMenn.Rows.Item(19).Insert()

# Insert a blank row before the (old) row 15 - this is where the new "Kule"
# entry belongs (at the very top of this block, pushing everything else,
# including the row just inserted above, one row further down).
$wsMenn.Rows.Item(15).Insert()

# New row 15: Marcel Grohmann's new best "Kule" result.
$wsMenn.Range("B15").Value = "Vestby IL"
$wsMenn.Range("C15").Value = "Akershus"
$wsMenn.Range("D15").Value = 303
$wsMenn.Range("E15").Value = "Kule"
$wsMenn.Range("F15").Value = "Marcel Grohmann"
$wsMenn.Range("G15").Value = 1977
$wsMenn.Range("H15").Value = "8,12"
$wsMenn.Range("I15").Value = "Bærum"
$wsMenn.Range("J15").Value = "29.10.2022"

# Row 16 (previously row 15, shifted down).
$wsMenn.Range("B16").Value = "Vestby IL"
$wsMenn.Range("C16").Value = "Akershus"
$wsMenn.Range("D16").Value = 294
$wsMenn.Range("E16").Value = "Kule"
$wsMenn.Range("F16").Value = "Marcel Grohmann"
$wsMenn.Range("G16").Value = 1977
$wsMenn.Range("H16").Value = "8,01"
$wsMenn.Range("I16").Value = "Ski"
# Format as text first so the day<=12 date-like string "01.10.2022" is not
# auto-converted into a date serial number by Excel's smart input parsing.
$wsMenn.Range("J16").NumberFormat = "@"
$wsMenn.Range("J16").Value = "01.10.2022"

# Row 17 (previously row 16, shifted down).
$wsMenn.Range("B17").Value = "Vestby IL"
$wsMenn.Range("C17").Value = "Akershus"
$wsMenn.Range("D17").Value = 221
$wsMenn.Range("E17").Value = "Spyd"
$wsMenn.Range("F17").Value = "Marcel Grohmann"
$wsMenn.Range("G17").Value = 1977
$wsMenn.Range("H17").Value = "23,90"
$wsMenn.Range("I17").Value = "Ski"
$wsMenn.Range("J17").NumberFormat = "@"
$wsMenn.Range("J17").Value = "01.10.2022"

# Row 18 (previously row 17, shifted down).
$wsMenn.Range("B18").Value = "Skogsvåg IL"
$wsMenn.Range("C18").Value = "Hordaland"
$wsMenn.Range("D18").Value = 190
$wsMenn.Range("E18").Value = "Lengde u.t"
$wsMenn.Range("F18").Value = "Kenneth Sangolt"
$wsMenn.Range("G18").Value = 1981
$wsMenn.Range("H18").Value = "2,31"
$wsMenn.Range("I18").Value = "Leikvang"
$wsMenn.Range("J18").Value = "19.03.2022"

# Row 19 (previously row 18, shifted down).
$wsMenn.Range("B19").Value = "Jægervatnet IL"
$wsMenn.Range("C19").Value = "Troms"
$wsMenn.Range("D19").Value = 180
$wsMenn.Range("E19").Value = "Høyde u.t"
$wsMenn.Range("F19").Value = "Kjell Ivar Robertsen"
$wsMenn.Range("G19").Value = 1964
$wsMenn.Range("H19").Value = "1,15"
$wsMenn.Range("I19").Value = "Grimstad"
$wsMenn.Range("J19").Value = "26.03.2022"

# New row 20: Marcel Grohmann's new "Diskos" result.
$wsMenn.Range("B20").Value = "Vestby IL"
$wsMenn.Range("C20").Value = "Akershus"
$wsMenn.Range("D20").Value = 114
$wsMenn.Range("E20").Value = "Diskos"
$wsMenn.Range("F20").Value = "Marcel Grohmann"
$wsMenn.Range("G20").Value = 1977
$wsMenn.Range("H20").Value = "16,63"
$wsMenn.Range("I20").Value = "Bærum"
$wsMenn.Range("J20").Value = "29.10.2022"

# Row 21 (previously row 19, shifted down).
$wsMenn.Range("B21").Value = "Samhald IL"
$wsMenn.Range("C21").Value = "Møre og Romsdal"
$wsMenn.Range("D21").Value = 34
$wsMenn.Range("E21").Value = "Lengde"
$wsMenn.Range("F21").Value = "Oddbjørn Bergheim"
$wsMenn.Range("G21").Value = 1950
$wsMenn.Range("H21").Value = "3,68"
$wsMenn.Range("I21").Value = "Grimstad"
$wsMenn.Range("J21").Value = "26.03.2022"

# Row 22 (previously row 20, shifted down).
$wsMenn.Range("B22").Value = "Samhald IL"
$wsMenn.Range("C22").Value = "Møre og Romsdal"
$wsMenn.Range("D22").Value = 32
$wsMenn.Range("E22").Value = "Tresteg"
$wsMenn.Range("F22").Value = "Oddbjørn Bergheim"
$wsMenn.Range("G22").Value = 1950
$wsMenn.Range("H22").Value = "7,60"
$wsMenn.Range("I22").Value = "Tromsø"
$wsMenn.Range("J22").Value = "21.08.2022"

# Row 23 (previously row 21, shifted down).
$wsMenn.Range("B23").Value = "Jægervatnet IL"
$wsMenn.Range("C23").Value = "Troms"
$wsMenn.Range("D23").Value = 26
$wsMenn.Range("E23").Value = "Lengde"
$wsMenn.Range("F23").Value = "Kjell Ivar Robertsen"
$wsMenn.Range("G23").Value = 1964
$wsMenn.Range("H23").Value = "3,65"
$wsMenn.Range("I23").Value = "Grimstad"
$wsMenn.Range("J23").Value = "26.03.2022"

# Row 24 (previously row 22, shifted down).
$wsMenn.Range("B24").Value = "Jægervatnet IL"
$wsMenn.Range("C24").Value = "Troms"
$wsMenn.Range("D24").Value = 10
$wsMenn.Range("E24").Value = "Lengde u.t"
$wsMenn.Range("F24").Value = "Kjell Ivar Robertsen"
$wsMenn.Range("G24").Value = 1964
$wsMenn.Range("H24").Value = "2,13"
$wsMenn.Range("I24").Value = "Grimstad"
$wsMenn.Range("J24").Value = "26.03.2022"

# ---------------------------------------------------------------------------
# Sheet "klubbres": update Marcel Grohmann's "Kule" result to the new value,
# and add his new "Diskos" result as a new row at the end of his club's
# results table.
# ---------------------------------------------------------------------------
$wsKlubb = $wb.Worksheets.Item("klubbres")

# Update existing "Kule" row with new result.
$wsKlubb.Range("E63").Value = "8,12"
$wsKlubb.Range("F63").Value = 303
$wsKlubb.Range("G63").Value = "Bærum"
$wsKlubb.Range("H63").Value = "29.10.2022"

# New row 65: Marcel Grohmann's new "Diskos" result.
$wsKlubb.Range("B65").Value = "Marcel Grohmann"
$wsKlubb.Range("C65").Value = 1977
$wsKlubb.Range("D65").Value = "Diskos"
$wsKlubb.Range("E65").Value = "16,63"
$wsKlubb.Range("F65").Value = 114
$wsKlubb.Range("G65").Value = "Bærum"
$wsKlubb.Range("H65").Value = "29.10.2022"
